# Updated cryptos list - applies price/volume/coin changes per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.476.64"
$ws.Range("D3").Value = "1.951.74"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.60"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4775"
$ws.Range("E7").Value = "  -4.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4022"
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.60"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08499"
$ws.Range("E10").Value = "  -5.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.058"
$ws.Range("E11").Value = "  -5.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.07"
$ws.Range("E12").Value = "  -5.13%  "
$ws.Range("D13").Value = "1.979.03"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.616"
$ws.Range("E14").Value = "  -5.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.189"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.014"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.04"
$ws.Range("E18").Value = "  -5.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06625"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.75"
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.809"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "28.521.14"
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.52"
$ws.Range("E24").Value = "  -3.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "2.189.67"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.12"
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.951"
$ws.Range("E29").Value = "  -7.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.160"
$ws.Range("E30").Value = "  -6.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.72"
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9951"
$ws.Range("E32").Value = "  -5.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09557"
$ws.Range("E33").Value = "  -4.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.449"
$ws.Range("E34").Value = "  -7.51%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.601"
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.663"
$ws.Range("E36").Value = "  -3.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02341"
$ws.Range("E37").Value = "  -5.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06228"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.770"
$ws.Range("E39").Value = "  -5.35%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.261"
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6233"
$ws.Range("E41").Value = "  -4.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.11"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.011"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1928"
$ws.Range("E44").Value = "  -6.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.328"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5971"
$ws.Range("E46").Value = "  -6.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.97"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.064"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.409"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06811"
$ws.Range("E51").Value = "  -2.54%  "
